$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix umlaut: "neue Straße" -> "neue Strasse"
$ws.Range("B4").Value = "neue Strasse"

# Move active selection to B5 (reflects user interaction captured in the diff)
$ws.Range("B5").Select()
